# Applies the row swaps/rotations for the Bosnia Herzegovina Premier Liga sheet
# as described by the commit diff (rows 36/37, 76/77, 122/123, 189-194, 197/198).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 (id stays 34)
$ws.Range("B36").Value = 6865299
$ws.Range("E36").Value = "Siroki Brijeg"
$ws.Range("F36").Value = "Zvijezda 09"
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = "H"
$ws.Range("L36").Value = 1.25
$ws.Range("M36").Value = 5.5
$ws.Range("N36").Value = 8
$ws.Range("O36").Value = 1.4
$ws.Range("P36").Value = 4.75
$ws.Range("Q36").Value = 5.75
$ws.Range("R36").Value = -1.25
$ws.Range("S36").Value = 1.9
$ws.Range("T36").Value = 1.9
$ws.Range("U36").Value = 2.75
$ws.Range("V36").Value = 1.85
$ws.Range("W36").Value = 1.95
$ws.Range("X36").Value = 0.3999999999999999
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = -1
$ws.Range("AA36").Value = -0.5
$ws.Range("AB36").Value = 0.45
$ws.Range("AC36").Value = 0.425
$ws.Range("AD36").Value = -0.5

# Row 37 (id stays 35)
$ws.Range("B37").Value = 6864629
$ws.Range("E37").Value = "Borac Banja Luka"
$ws.Range("F37").Value = "NK Posusje"
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = "H"
$ws.Range("L37").Value = 1.363
$ws.Range("M37").Value = 4.5
$ws.Range("N37").Value = 6.5
$ws.Range("O37").Value = 1.363
$ws.Range("P37").Value = 4.2
$ws.Range("Q37").Value = 6.5
$ws.Range("R37").Value = -1.25
$ws.Range("S37").Value = 1.95
$ws.Range("T37").Value = 1.85
$ws.Range("U37").Value = 2.5
$ws.Range("V37").Value = 1.925
$ws.Range("W37").Value = 1.875
$ws.Range("X37").Value = 0.363
$ws.Range("Y37").Value = -1
$ws.Range("Z37").Value = -1
$ws.Range("AA37").Value = -0.5
$ws.Range("AB37").Value = 0.425
$ws.Range("AC37").Value = -1
$ws.Range("AD37").Value = 0.875

# Row 76 (id stays 74)
$ws.Range("B76").Value = 6865328
$ws.Range("E76").Value = "Siroki Brijeg"
$ws.Range("F76").Value = "NK Posusje"
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1
$ws.Range("K76").Value = "D"
$ws.Range("L76").Value = 2
$ws.Range("M76").Value = 3
$ws.Range("N76").Value = 3.5
$ws.Range("O76").Value = 2.1
$ws.Range("P76").Value = 3
$ws.Range("Q76").Value = 3.3
$ws.Range("R76").Value = -0.25
$ws.Range("S76").Value = 1.825
$ws.Range("T76").Value = 1.975
$ws.Range("U76").Value = 2
$ws.Range("V76").Value = 1.825
$ws.Range("W76").Value = 1.975
$ws.Range("X76").Value = -1
$ws.Range("Y76").Value = 2
$ws.Range("Z76").Value = -1
$ws.Range("AA76").Value = -0.5
$ws.Range("AB76").Value = 0.4875
$ws.Range("AC76").Value = 0
$ws.Range("AD76").Value = 0

# Row 77 (id stays 75)
$ws.Range("B77").Value = 6865377
$ws.Range("E77").Value = "Zrinjski Mostar"
$ws.Range("F77").Value = "FK Tuzla City"
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 1
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = "H"
$ws.Range("L77").Value = 1.333
$ws.Range("M77").Value = 5
$ws.Range("N77").Value = 6
$ws.Range("O77").Value = 1.166
$ws.Range("P77").Value = 6.5
$ws.Range("Q77").Value = 13
$ws.Range("R77").Value = -2
$ws.Range("S77").Value = 1.9
$ws.Range("T77").Value = 1.9
$ws.Range("U77").Value = 3.25
$ws.Range("V77").Value = 1.95
$ws.Range("W77").Value = 1.85
$ws.Range("X77").Value = 0.1659999999999999
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = -1
$ws.Range("AA77").Value = 0
$ws.Range("AB77").Value = 0
$ws.Range("AC77").Value = 0.95
$ws.Range("AD77").Value = -1

# Row 122 (id stays 120)
$ws.Range("B122").Value = 6865381
$ws.Range("E122").Value = "FK Tuzla City"
$ws.Range("F122").Value = "Zvijezda 09"
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = "H"
$ws.Range("L122").Value = 1.666
$ws.Range("M122").Value = 3.6
$ws.Range("N122").Value = 4.333
$ws.Range("O122").Value = 1.5
$ws.Range("P122").Value = 4
$ws.Range("Q122").Value = 5.25
$ws.Range("R122").Value = -1
$ws.Range("S122").Value = 1.925
$ws.Range("T122").Value = 1.875
$ws.Range("U122").Value = 2.5
$ws.Range("V122").Value = 1.8
$ws.Range("W122").Value = 2
$ws.Range("X122").Value = 0.5
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = -1
$ws.Range("AA122").Value = 0.925
$ws.Range("AB122").Value = -1
$ws.Range("AC122").Value = -1
$ws.Range("AD122").Value = 1

# Row 123 (id stays 121)
$ws.Range("B123").Value = 6865363
$ws.Range("E123").Value = "NK Igman Konjic"
$ws.Range("F123").Value = "Siroki Brijeg"
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = "H"
$ws.Range("L123").Value = 2
$ws.Range("M123").Value = 3.3
$ws.Range("N123").Value = 3.25
$ws.Range("O123").Value = 2.3
$ws.Range("P123").Value = 3.2
$ws.Range("Q123").Value = 2.75
$ws.Range("R123").Value = -0.25
$ws.Range("S123").Value = 2.05
$ws.Range("T123").Value = 1.75
$ws.Range("U123").Value = 2
$ws.Range("V123").Value = 1.9
$ws.Range("W123").Value = 1.9
$ws.Range("X123").Value = 1.3
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = -1
$ws.Range("AA123").Value = 1.05
$ws.Range("AB123").Value = -1
$ws.Range("AC123").Value = -1
$ws.Range("AD123").Value = 0.8999999999999999

# Row 189 (id stays 187)
$ws.Range("B189").Value = 7952777
$ws.Range("E189").Value = "Borac Banja Luka"
$ws.Range("F189").Value = "NK Igman Konjic"
$ws.Range("G189").Value = 4
$ws.Range("H189").Value = 3
$ws.Range("I189").Value = 1
$ws.Range("J189").Value = 2
$ws.Range("K189").Value = "H"
$ws.Range("L189").Value = 1.25
$ws.Range("M189").Value = 5.75
$ws.Range("N189").Value = 7
$ws.Range("O189").Value = 1.2
$ws.Range("P189").Value = 5.75
$ws.Range("Q189").Value = 12
$ws.Range("R189").Value = -2
$ws.Range("S189").Value = 1.95
$ws.Range("T189").Value = 1.85
$ws.Range("U189").Value = 3.25
$ws.Range("V189").Value = 1.9
$ws.Range("W189").Value = 1.9
$ws.Range("X189").Value = 0.2
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = -1
$ws.Range("AA189").Value = -1
$ws.Range("AB189").Value = 0.8500000000000001
$ws.Range("AC189").Value = 0.8999999999999999
$ws.Range("AD189").Value = -1

# Row 190 (id stays 188)
$ws.Range("B190").Value = 7952780
$ws.Range("E190").Value = "Velez Mostar"
$ws.Range("F190").Value = "GOSK Gabela"
$ws.Range("G190").Value = 3
$ws.Range("H190").Value = 3
$ws.Range("I190").Value = 1
$ws.Range("J190").Value = 1
$ws.Range("K190").Value = "D"
$ws.Range("L190").Value = 1.4
$ws.Range("M190").Value = 4
$ws.Range("N190").Value = 7
$ws.Range("O190").Value = 1.363
$ws.Range("P190").Value = 4.2
$ws.Range("Q190").Value = 8
$ws.Range("R190").Value = -1.5
$ws.Range("S190").Value = 2
$ws.Range("T190").Value = 1.8
$ws.Range("U190").Value = 2.75
$ws.Range("V190").Value = 1.825
$ws.Range("W190").Value = 1.975
$ws.Range("X190").Value = -1
$ws.Range("Y190").Value = 3.2
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = 0.8
$ws.Range("AC190").Value = 0.825
$ws.Range("AD190").Value = -1

# Row 191 (id stays 189)
$ws.Range("B191").Value = 7952779
$ws.Range("E191").Value = "Zrinjski Mostar"
$ws.Range("F191").Value = "FK Tuzla City"
$ws.Range("G191").Value = 4
$ws.Range("H191").Value = 0
$ws.Range("I191").Value = 2
$ws.Range("J191").Value = 0
$ws.Range("K191").Value = "H"
$ws.Range("L191").Value = 1.25
$ws.Range("M191").Value = 5.75
$ws.Range("N191").Value = 7
$ws.Range("O191").Value = 1.055
$ws.Range("P191").Value = 13
$ws.Range("Q191").Value = 17
$ws.Range("R191").Value = -3.5
$ws.Range("S191").Value = 1.975
$ws.Range("T191").Value = 1.825
$ws.Range("U191").Value = 4.75
$ws.Range("V191").Value = 1.825
$ws.Range("W191").Value = 1.975
$ws.Range("X191").Value = 0.05499999999999994
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = 0.9750000000000001
$ws.Range("AB191").Value = -1
$ws.Range("AC191").Value = -1
$ws.Range("AD191").Value = 0.9750000000000001

# Row 192 (id stays 190)
$ws.Range("B192").Value = 7952778
$ws.Range("E192").Value = "Sloga"
$ws.Range("F192").Value = "Siroki Brijeg"
$ws.Range("G192").Value = 2
$ws.Range("H192").Value = 3
$ws.Range("I192").Value = 2
$ws.Range("J192").Value = 2
$ws.Range("K192").Value = "A"
$ws.Range("L192").Value = 1.727
$ws.Range("M192").Value = 3.75
$ws.Range("N192").Value = 3.75
$ws.Range("O192").Value = 1.7
$ws.Range("P192").Value = 3.9
$ws.Range("Q192").Value = 3.9
$ws.Range("R192").Value = -0.75
$ws.Range("S192").Value = 1.975
$ws.Range("T192").Value = 1.825
$ws.Range("U192").Value = 2.25
$ws.Range("V192").Value = 1.8
$ws.Range("W192").Value = 2
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = -1
$ws.Range("Z192").Value = 2.9
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = 0.825
$ws.Range("AC192").Value = 0.8
$ws.Range("AD192").Value = -1

# Row 193 (id stays 191)
$ws.Range("B193").Value = 7952781
$ws.Range("E193").Value = "Zvijezda 09"
$ws.Range("F193").Value = "Zeljeznicar"
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 5
$ws.Range("I193").Value = 0
$ws.Range("J193").Value = 1
$ws.Range("K193").Value = "A"
$ws.Range("L193").Value = 2.15
$ws.Range("M193").Value = 3.25
$ws.Range("N193").Value = 2.9
$ws.Range("O193").Value = 3.6
$ws.Range("P193").Value = 3.4
$ws.Range("Q193").Value = 1.85
$ws.Range("R193").Value = 0.5
$ws.Range("S193").Value = 1.875
$ws.Range("T193").Value = 1.925
$ws.Range("U193").Value = 2.5
$ws.Range("V193").Value = 1.975
$ws.Range("W193").Value = 1.825
$ws.Range("X193").Value = -1
$ws.Range("Y193").Value = -1
$ws.Range("Z193").Value = 0.8500000000000001
$ws.Range("AA193").Value = -1
$ws.Range("AB193").Value = 0.925
$ws.Range("AC193").Value = 0.9750000000000001
$ws.Range("AD193").Value = -1

# Row 194 (id stays 192)
$ws.Range("B194").Value = 7952776
$ws.Range("E194").Value = "FK Sarajevo"
$ws.Range("F194").Value = "NK Posusje"
$ws.Range("G194").Value = 1
$ws.Range("H194").Value = 1
$ws.Range("I194").Value = 0
$ws.Range("J194").Value = 0
$ws.Range("K194").Value = "D"
$ws.Range("L194").Value = 1.571
$ws.Range("M194").Value = 3.4
$ws.Range("N194").Value = 5.5
$ws.Range("O194").Value = 1.363
$ws.Range("P194").Value = 3.9
$ws.Range("Q194").Value = 8
$ws.Range("R194").Value = -1.25
$ws.Range("S194").Value = 1.85
$ws.Range("T194").Value = 1.95
$ws.Range("U194").Value = 2.75
$ws.Range("V194").Value = 1.925
$ws.Range("W194").Value = 1.875
$ws.Range("X194").Value = -1
$ws.Range("Y194").Value = 2.9
$ws.Range("Z194").Value = -1
$ws.Range("AA194").Value = -1
$ws.Range("AB194").Value = 0.95
$ws.Range("AC194").Value = -1
$ws.Range("AD194").Value = 0.875

# Row 197 (id stays 195)
$ws.Range("B197").Value = 8259814
$ws.Range("E197").Value = "Siroki Brijeg"
$ws.Range("F197").Value = "FK Sarajevo"
$ws.Range("G197").Value = 2
$ws.Range("H197").Value = 2
$ws.Range("I197").Value = $null
$ws.Range("J197").Value = $null
$ws.Range("K197").Value = "D"
$ws.Range("L197").Value = 3.4
$ws.Range("M197").Value = 3.1
$ws.Range("N197").Value = 2
$ws.Range("O197").Value = 9.5
$ws.Range("P197").Value = 4.75
$ws.Range("Q197").Value = 1.25
$ws.Range("R197").Value = 1.5
$ws.Range("S197").Value = 1.975
$ws.Range("T197").Value = 1.825
$ws.Range("U197").Value = 2.75
$ws.Range("V197").Value = 1.95
$ws.Range("W197").Value = 1.85
$ws.Range("X197").Value = -1
$ws.Range("Y197").Value = 3.75
$ws.Range("Z197").Value = -1
$ws.Range("AA197").Value = 0.9750000000000001
$ws.Range("AB197").Value = -1
$ws.Range("AC197").Value = 0.95
$ws.Range("AD197").Value = -1

# Row 198 (id stays 196)
$ws.Range("B198").Value = 8259815
$ws.Range("E198").Value = "NK Posusje"
$ws.Range("F198").Value = "Zvijezda 09"
$ws.Range("G198").Value = 2
$ws.Range("H198").Value = 0
$ws.Range("I198").Value = $null
$ws.Range("J198").Value = $null
$ws.Range("K198").Value = "H"
$ws.Range("L198").Value = 1.4
$ws.Range("M198").Value = 4
$ws.Range("N198").Value = 6.5
$ws.Range("O198").Value = 1.25
$ws.Range("P198").Value = 5
$ws.Range("Q198").Value = 8.5
$ws.Range("R198").Value = -1.75
$ws.Range("S198").Value = 2
$ws.Range("T198").Value = 1.8
$ws.Range("U198").Value = 3
$ws.Range("V198").Value = 1.925
$ws.Range("W198").Value = 1.875
$ws.Range("X198").Value = 0.25
$ws.Range("Y198").Value = -1
$ws.Range("Z198").Value = -1
$ws.Range("AA198").Value = 0.5
$ws.Range("AB198").Value = -0.5
$ws.Range("AC198").Value = -1
$ws.Range("AD198").Value = 0.875
